# Update the marksheet's correct/total marks values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row: correct answer score changed from 3 to 5
$ws.Range("B11").Value = 5

# Total row: total marks changed from 48 to 80
$ws.Range("B12").Value = 80

# Total row: Corr/Total marks text changed from "44/84" to "80/140"
$ws.Range("E12").Value = "80/140"

$wb.Save()
